# Apply updated cryptos list values (price + 1h volume change) per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.245.29"
$ws.Range("E2").Value = "  -0.09%  "

$ws.Range("D3").Value = "1.591.98"
$ws.Range("E3").Value = "  +0.13%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'212.62"
$ws.Range("E5").Value = "  -0.19%  "

$ws.Range("E6").Value = "  -0.30%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  -0.52%  "

$ws.Range("D9").Value = "'0.0606"
$ws.Range("E9").Value = "  -0.54%  "

$ws.Range("D10").Value = "'18.92"
$ws.Range("E10").Value = "  -2.19%  "

$ws.Range("E11").Value = "  +0.24%  "

$ws.Range("D12").Value = "1.815.63"
$ws.Range("E12").Value = "  +0.12%  "

$ws.Range("D13").Value = "1.598.92"
$ws.Range("E13").Value = "  +0.44%  "

$ws.Range("E14").Value = "  -1.20%  "

$ws.Range("E15").Value = "  -2.72%  "

$ws.Range("D16").Value = "'63.86"
$ws.Range("E16").Value = "  -0.93%  "

$ws.Range("D17").Value = "26.253.44"
$ws.Range("E17").Value = "  -0.10%  "

$ws.Range("D18").Value = "0.0₃0722"
$ws.Range("E18").Value = "  -0.61%  "

$ws.Range("D19").Value = "'214.81"
$ws.Range("E19").Value = "  +0.60%  "

$ws.Range("D20").Value = "'7.37"
$ws.Range("E20").Value = "  -1.44%  "

$ws.Range("E21").Value = "  +0.09%  "

$ws.Range("E23").Value = "  +0.30%  "

$ws.Range("E24").Value = "  -2.68%  "

$ws.Range("D25").Value = "'144.70"
$ws.Range("E25").Value = "  -0.22%  "

$ws.Range("E26").Value = "  +0.05%  "

$ws.Range("E27").Value = "  -1.31%  "

$ws.Range("D29").Value = "'15.11"
$ws.Range("E29").Value = "  -0.56%  "

$ws.Range("D30").Value = "'0.0494"
$ws.Range("E30").Value = "  -1.08%  "

$ws.Range("E31").Value = "  -0.16%  "

$ws.Range("D32").Value = "'3.20"
$ws.Range("E32").Value = "  -0.42%  "

$ws.Range("D33").Value = "1.418.45"
$ws.Range("E33").Value = "  +5.91%  "

$ws.Range("E34").Value = "  -0.24%  "

$ws.Range("E35").Value = "  -0.92%  "

$ws.Range("E36").Value = "  -1.34%  "

$ws.Range("E37").Value = "  -4.53%  "

$ws.Range("E38").Value = "  -0.65%  "

$ws.Range("D39").Value = "'0.824"
$ws.Range("E39").Value = "  +0.94%  "

$ws.Range("E40").Value = "  -0.12%  "

$ws.Range("E41").Value = "  +0.06%  "

$ws.Range("E42").Value = "  +0.99%  "

$ws.Range("D43").Value = "'0.937"
$ws.Range("E43").Value = "  -7.69%  "

$ws.Range("D44").Value = "'0.761"
$ws.Range("E44").Value = "  -0.21%  "

$ws.Range("D45").Value = "1.728.70"
$ws.Range("E45").Value = "  +0.24%  "

$ws.Range("D46").Value = "'60.72"
$ws.Range("E46").Value = "  -1.70%  "

$ws.Range("D47").Value = "'86.62"
$ws.Range("E47").Value = "  -0.82%  "

$ws.Range("E48").Value = "  -1.55%  "

$ws.Range("E49").Value = "  -0.68%  "

$ws.Range("D50").Value = "'0.0951"
$ws.Range("E50").Value = "  -2.88%  "

$ws.Range("D51").Value = "'0.999"
$ws.Range("E51").Value = "  -0.13%  "
